$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "参数" column (column A), shifting 报警/结果/参考范围 left.
$ws.Columns("A").Delete()

# The "结果" column (now column B) holds numeric-looking readings that must
# stay text (units were stripped from them), so pre-format as Text before
# writing the values - otherwise Excel auto-converts "85.1" etc. to numbers.
# (Row 9's "H17.1" is already non-numeric text and needs no such coercion.)
$ws.Range("B2:B8").NumberFormat = "@"
$ws.Range("B10").NumberFormat = "@"

# Row 1 - headers (already correct after shift, but set explicitly to be safe)
$ws.Range("A1").Value = "报警"
$ws.Range("B1").Value = "结果"
$ws.Range("C1").Value = "参考范围"

# Row 2
$ws.Range("A2").Value = ""
$ws.Range("B2").Value = "85.1"
$ws.Range("C2").Value = "82.0-95.0"

# Row 3
$ws.Range("A3").Value = ""
$ws.Range("B3").Value = "27.0"
$ws.Range("C3").Value = "27.0-31.0"

# Row 4
$ws.Range("A4").Value = ""
$ws.Range("B4").Value = "318"
$ws.Range("C4").Value = "320-360"

# Row 5
$ws.Range("A5").Value = ""
$ws.Range("B5").Value = "13.2"
$ws.Range("C5").Value = "11.5-14.5"

# Row 6
$ws.Range("A6").Value = ""
$ws.Range("B6").Value = "42.8"
$ws.Range("C6").Value = "35.0-56.0"

# Row 7
$ws.Range("A7").Value = ""
$ws.Range("B7").Value = "220"
$ws.Range("C7").Value = "100-300"

# Row 8
$ws.Range("A8").Value = ""
$ws.Range("B8").Value = "8.5"
$ws.Range("C8").Value = "7.0-11.0"

# Row 9
$ws.Range("A9").Value = "H17.1"
$ws.Range("B9").Value = "H17.1"
$ws.Range("C9").Value = "15.0-17.0"

# Row 10
$ws.Range("A10").Value = ""
$ws.Range("B10").Value = "0.187"
$ws.Range("C10").Value = "0.108-0.282"
